# Add a new worksheet named "Sheet1" after the existing "STM32F4" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("STM32F4")

# Add the new data row to STM32F4 first (while it's still the active sheet)
$ws1.Range("A29").Value = "le_Time::Decode()"
$ws1.Range("B29").Value = "8 frames"
$ws1.Range("E29").Value = 3830
$ws1.Range("F29").Value = "exc_time = 45.6 us (max)"

# Change the selection on STM32F4 to A30 (still on that sheet, before switching away)
$ws1.Range("A30").Select()

# Insert the new worksheet after STM32F4; this becomes the active sheet/tab
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Sheet1"
